$d = $word.ActiveDocument

# First paragraph of the document: the hidden bookmark/ID placeholder line.
$p1 = $d.Paragraphs(1)
$start = $p1.Range.Start
$end = $p1.Range.End

# Add a (space-only, no line) paragraph border matching the rest of the
# document's body paragraphs, and widen the left indent to match them too.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# Replace the two runs ("**ID__AFFARS_pgi_5315_topic_9__ID**" + " ") with a
# single run containing the updated placeholder text (keeps the first run's
# character formatting, drops the trailing space-only run).
$r1 = $d.Range($start, $end)
$r1.Text = "**ID__AFFARS_AF_PGI_5315_407_91__ID**"
